$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.05
$ws.Range("O2").Value = 1.33
$ws.Range("U2").Value = 1.92
$ws.Range("V2").Value = 1.77
$ws.Range("M3").Value = 1.07
$ws.Range("O3").Value = 1.37
$ws.Range("P3").Value = 3
$ws.Range("U3").Value = 1.92
$ws.Range("V3").Value = 1.77
$ws.Range("BD4").Value = 151
$ws.Range("G5").Value = 2.6
$ws.Range("I5").Value = 3.1
$ws.Range("J5").Value = 3.4
$ws.Range("M5").Value = 1.13
$ws.Range("N5").Value = 6
$ws.Range("Z5").Value = 26
$ws.Range("AA5").Value = 26
$ws.Range("AC5").Value = 6
$ws.Range("AI5").Value = 13
$ws.Range("AJ5").Value = 12
$ws.Range("AS5").Value = 301
$ws.Range("G7").Value = 2.6
$ws.Range("I7").Value = 2.6
$ws.Range("K7").Value = 2.1
$ws.Range("L7").Value = 3.25
$ws.Range("Q7").Value = 1.95
$ws.Range("R7").Value = 1.9
$ws.Range("U7").Value = 1.75
$ws.Range("V7").Value = 2
$ws.Range("AZ7").Value = 51
$ws.Range("O8").Value = 1.25
$ws.Range("P8").Value = 3.75
$ws.Range("Q8").Value = 1.88
$ws.Range("R8").Value = 1.98
$ws.Range("U8").Value = 1.77
$ws.Range("V8").Value = 1.87
$ws.Range("Q9").Value = 1.72
$ws.Range("G11").Value = 1.67
$ws.Range("H11").Value = 4.2
$ws.Range("J11").Value = 2.2
$ws.Range("K11").Value = 2.5
$ws.Range("Q11").Value = 1.5
$ws.Range("R11").Value = 2.4
$ws.Range("S11").Value = 1.25
$ws.Range("T11").Value = 3.75
$ws.Range("W11").Value = 10
$ws.Range("X11").Value = 10
$ws.Range("AB11").Value = 19
$ws.Range("AD11").Value = 8.5
$ws.Range("AL11").Value = 29
$ws.Range("AM11").Value = 29
$ws.Range("AP11").Value = 15
$ws.Range("AQ11").Value = 23
$ws.Range("AS11").Value = 81
$ws.Range("AT11").Value = 3.75
$ws.Range("AY11").Value = 23
$ws.Range("BA11").Value = 67
$ws.Range("BB11").Value = 126
$ws.Range("Q12").Value = 1.63
$ws.Range("Q15").Value = 1.65
$ws.Range("R15").Value = 2.2
$ws.Range("G16").Value = 1.53
$ws.Range("I16").Value = 5.5
$ws.Range("J16").Value = 2.1
$ws.Range("AX16").Value = 29
$ws.Range("G17").Value = 2.4
$ws.Range("I17").Value = 2.55
$ws.Range("J17").Value = 2.88
$ws.Range("L17").Value = 3
$ws.Range("O17").Value = 1.11
$ws.Range("P17").Value = 6.5
$ws.Range("Q17").Value = 1.4
$ws.Range("R17").Value = 2.88
$ws.Range("AC17").Value = 23
$ws.Range("AJ17").Value = 11
$ws.Range("AK17").Value = 29
$ws.Range("AO17").Value = 12
$ws.Range("AX17").Value = 13
$ws.Range("BC17").Value = 201
$ws.Range("U18").Value = 1.69
$ws.Range("Q19").Value = 1.63
$ws.Range("U19").Value = 1.63
$ws.Range("I20").Value = 1.44
$ws.Range("Q20").Value = 1.44
$ws.Range("U20").Value = 1.63
$ws.Range("G21").Value = 1.36
$ws.Range("Q21").Value = 1.3
$ws.Range("U21").Value = 1.5
$ws.Range("V21").Value = 2.37
$ws.Range("G22").Value = 2.45
$ws.Range("I22").Value = 2.63
$ws.Range("K22").Value = 2.25
$ws.Range("L22").Value = 3.2
$ws.Range("Q22").Value = 1.67
$ws.Range("U22").Value = 1.54
$ws.Range("W22").Value = 10
$ws.Range("X22").Value = 13
$ws.Range("Z22").Value = 23
$ws.Range("AG22").Value = 151
$ws.Range("AL22").Value = 21
$ws.Range("AM22").Value = 26
$ws.Range("AX22").Value = 15
$ws.Range("Q26").Value = 1.75
$ws.Range("R26").Value = 2.05
$ws.Range("J32").Value = 2.88
$ws.Range("K32").Value = 2.38
$ws.Range("R33").Value = 1.58
$ws.Range("G34").Value = 1.57
$ws.Range("Q34").Value = 1.77
$ws.Range("R34").Value = 1.97
$ws.Range("G35").Value = 2.55
$ws.Range("I35").Value = 2.9
$ws.Range("J35").Value = 3.25
$ws.Range("L35").Value = 3.6
$ws.Range("R35").Value = 1.62
$ws.Range("X35").Value = 12
$ws.Range("AI35").Value = 13
$ws.Range("AJ35").Value = 11
$ws.Range("AK35").Value = 29
$ws.Range("AM35").Value = 34
$ws.Range("AO35").Value = 15
$ws.Range("AQ35").Value = 51
$ws.Range("AW35").Value = 4.75
